$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that looks numeric (e.g. "1.009", "46.80").
# Force it to Text format first so Excel does not coerce these into
# floating point numbers and drop the significant trailing/format digits.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "27.768.45"
$ws.Range("E2").Value = "  +0.94%  "

# Row 3
$ws.Range("D3").Value = "1.847.37"
$ws.Range("E3").Value = "  +0.20%  "

# Row 4
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").Value = "335.53"
$ws.Range("E5").Value = "  +0.44%  "

# Row 6
$ws.Range("D6").Value = "1.008"
$ws.Range("E6").Value = "  +0.08%  "

# Row 7
$ws.Range("D7").Value = "0.4649"
$ws.Range("E7").Value = "  +0.72%  "

# Row 8
$ws.Range("D8").Value = "0.3859"
$ws.Range("E8").Value = "  -0.14%  "

# Row 9
$ws.Range("D9").Value = "46.80"
$ws.Range("E9").Value = "  +1.80%  "

# Row 10
$ws.Range("D10").Value = "0.07892"
$ws.Range("E10").Value = "  -0.16%  "

# Row 11
$ws.Range("D11").Value = "0.9669"
$ws.Range("E11").Value = "  -3.22%  "

# Row 12
$ws.Range("D12").Value = "21.26"
$ws.Range("E12").Value = "  -1.15%  "

# Row 13
$ws.Range("D13").Value = "1.867.66"
$ws.Range("E13").Value = "  +0.88%  "

# Row 14
$ws.Range("D14").Value = "5.889"
$ws.Range("E14").Value = "  -1.33%  "

# Row 15
$ws.Range("D15").Value = "7.140"
$ws.Range("E15").Value = "  -0.12%  "

# Row 16
$ws.Range("D16").Value = "1.010"
$ws.Range("E16").Value = "  +0.06%  "

# Row 17
$ws.Range("D17").Value = "90.95"
$ws.Range("E17").Value = "  +3.01%  "

# Row 18
$ws.Range("D18").Value = "0.06601"
$ws.Range("E18").Value = "  -1.00%  "

# Row 19
$ws.Range("D19").Value = "0.00001028"
$ws.Range("E19").Value = "  -0.65%  "

# Row 20
$ws.Range("D20").Value = "17.27"
$ws.Range("E20").Value = "  +0.77%  "

# Row 21
$ws.Range("D21").Value = "1.008"
$ws.Range("E21").Value = "  +0.06%  "

# Row 22
$ws.Range("D22").Value = "27.765.98"
$ws.Range("E22").Value = "  +0.92%  "

# Row 23
$ws.Range("D23").Value = "5.353"
$ws.Range("E23").Value = "  -0.74%  "

# Row 24
$ws.Range("D24").Value = "10.81"
$ws.Range("E24").Value = "  -0.92%  "

# Row 25
$ws.Range("D25").Value = "2.300"
$ws.Range("E25").Value = "  -0.16%  "

# Row 26
$ws.Range("D26").Value = "2.086.59"
$ws.Range("E26").Value = "  +0.76%  "

# Row 27
$ws.Range("D27").Value = "158.86"
$ws.Range("E27").Value = "  -0.13%  "

# Row 28
$ws.Range("D28").Value = "19.42"
$ws.Range("E28").Value = "  -0.34%  "

# Row 29
$ws.Range("D29").Value = "2.065"
$ws.Range("E29").Value = "  -2.33%  "

# Row 30
$ws.Range("D30").Value = "5.372"
$ws.Range("E30").Value = "  -0.73%  "

# Row 31
$ws.Range("D31").Value = "118.67"
$ws.Range("E31").Value = "  -1.34%  "

# Row 32
$ws.Range("D32").Value = "0.09434"
$ws.Range("E32").Value = "  +0.39%  "

# Row 33
$ws.Range("D33").Value = "0.9410"
$ws.Range("E33").Value = "  -3.44%  "

# Row 34
$ws.Range("D34").Value = "3.602"
$ws.Range("E34").Value = "  +0.08%  "

# Row 35
$ws.Range("D35").Value = "5.253"
$ws.Range("E35").Value = "  -0.81%  "

# Row 36
$ws.Range("D36").Value = "1.325"
$ws.Range("E36").Value = "  -0.65%  "

# Row 37
$ws.Range("E37").Value = "  +0.21%  "

# Row 38
$ws.Range("D38").Value = "0.02214"
$ws.Range("E38").Value = "  -0.37%  "

# Row 39
$ws.Range("D39").Value = "8.223"
$ws.Range("E39").Value = "  -0.68%  "

# Row 40
$ws.Range("D40").Value = "1.007"
$ws.Range("E40").Value = "  +0.00%  "

# Row 41
$ws.Range("D41").Value = "1.151"
$ws.Range("E41").Value = "  -2.47%  "

# Row 42
$ws.Range("D42").Value = "0.5808"
$ws.Range("E42").Value = "  -1.63%  "

# Row 43
$ws.Range("D43").Value = "0.1846"
$ws.Range("E43").Value = "  -0.93%  "

# Row 44
$ws.Range("D44").Value = "10.06"
$ws.Range("E44").Value = "  -2.66%  "

# Row 45
$ws.Range("D45").Value = "1.283"
$ws.Range("E45").Value = "  +3.51%  "

# Row 46
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "11.99"
$ws.Range("E46").Value = "  -1.74%  "

# Row 47
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "0.5449"
$ws.Range("E47").Value = "  -2.42%  "

# Row 48
$ws.Range("D48").Value = "1.936"
$ws.Range("E48").Value = "  +1.43%  "

# Row 49
$ws.Range("D49").Value = "0.06853"
$ws.Range("E49").Value = "  +2.31%  "

# Row 50
$ws.Range("D50").Value = "110.70"
$ws.Range("E50").Value = "  +0.57%  "

# Row 51
$ws.Range("D51").Value = "1.008"
$ws.Range("E51").Value = "  -32.24%  "
